# Update the two date values in column B and move the active selection,
# as captured by the workbook diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B2: 2017-01-01 (42736) -> 2017-01-31 (42766)
$ws.Range("B2").Value = "1/31/2017"

# B3: 2017-01-02 (42737) -> 2017-02-28 (42794)
$ws.Range("B3").Value = "2/28/2017"

# Move the sheet's active cell/selection from A4 to H7
$ws.Range("H7").Select()
